$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule gains a new first task row ("Configuración tarjeta gráfica"),
# pushing the existing tasks down by one row, and every date becomes a
# plain stored value (no more "+15" formulas). Work from the bottom up so
# we never overwrite a source row before it has been copied from.

# 1) Push existing rows 3-6 down to rows 4-7, carrying their exact
#    formatting (border/number-format) along with them.
$ws.Cells.Item(6,1).Copy($ws.Cells.Item(7,1))
$ws.Cells.Item(6,2).Copy($ws.Cells.Item(7,2))

$ws.Cells.Item(5,1).Copy($ws.Cells.Item(6,1))
$ws.Cells.Item(5,2).Copy($ws.Cells.Item(6,2))

$ws.Cells.Item(4,1).Copy($ws.Cells.Item(5,1))
$ws.Cells.Item(4,2).Copy($ws.Cells.Item(5,2))

$ws.Cells.Item(3,1).Copy($ws.Cells.Item(4,1))
$ws.Cells.Item(3,2).Copy($ws.Cells.Item(4,2))

# 2) Row 3 keeps the same per-column formatting as the rest of the data
#    rows, so just reuse row 4 (its freshly-shifted copy) as the format
#    source for the new row.
$ws.Cells.Item(4,1).Copy($ws.Cells.Item(3,1))
$ws.Cells.Item(4,2).Copy($ws.Cells.Item(3,2))

# 3) Write the new task row.
$ws.Cells.Item(3,1).Value = "Configuración tarjeta gráfica"
$ws.Cells.Item(3,2).Value = 44093

# 4) Replace the old formula-driven dates with their plain stored values.
$ws.Cells.Item(4,2).Value = 44102
$ws.Cells.Item(5,2).Value = 44114
$ws.Cells.Item(6,2).Value = 44128
$ws.Cells.Item(7,2).Value = 44142

# 5) Re-confirm the task labels for the shifted rows (values already
#    carried over via Copy, but this keeps intent explicit/robust).
$ws.Cells.Item(4,1).Value = "Entrenamient"
$ws.Cells.Item(5,1).Value = "Pruebas"
$ws.Cells.Item(6,1).Value = "Evaluacion modelo"
$ws.Cells.Item(7,1).Value = "Analisis de resultados"

# 6) Column A widens slightly to fit the new, longer label.
$ws.Columns.Item(1).ColumnWidth = 25.25

# 7) Selection moves to B4 (no longer the merged title range).
$ws.Range("B4").Select()
